$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the note in A46: old text gets replaced with an extended version.
# (Excel automatically reindexes/cleans up the shared strings table on save,
# so editing just this cell reproduces the same shared-string shuffle seen
# in the target diff for A46:A49 and M53:N53.)
$ws.Cells.Item(46, 1).Value = "peak data is found in same way as above although may adjust parameters like approx_fsr if see fit"

# Narrow column L (12) down from its bestFit width to a plain custom width of 10.
# 9.14 "characters" is the COM ColumnWidth value that Excel stores internally as width=10.
$ws.Columns.Item(12).ColumnWidth = 9.14

# Add the new FWHM data run for sg_rr_20_025 2023-12-13 17-59-26 as row 54.
$ws.Cells.Item(54, 1).Value = "sg_rr_20_025 2023-12-13 17-59-26.csv"
$ws.Cells.Item(54, 2).Value = 0.01
$ws.Cells.Item(54, 3).Value = 1000
$ws.Cells.Item(54, 4).Value = 5001
$ws.Cells.Item(54, 5).Value = 1530
$ws.Cells.Item(54, 6).Value = 1570
$ws.Cells.Item(54, 7).Value = 0.5
$ws.Cells.Item(54, 8).Value = "(approx_fsr/2)/wavelength step size"
$ws.Cells.Item(54, 9).Value = 5
$ws.Cells.Item(54, 10).Value = 4.92
$ws.Cells.Item(54, 11).Value = 0.025354627641843101
$ws.Cells.Item(54, 12).Value = "yes"
$ws.Cells.Item(54, 13).Value = 0.14651376066498201
$ws.Cells.Item(54, 14).Value = 0.0162754097761134

# Update the view: scroll/select near the newly added row and zoom out a bit.
$ws.Activate()
$ws.Range("A30").Select()
$excel.ActiveWindow.Zoom = 76
